$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the numeric values in B2 and B3 with text values "A" and "B"
$ws.Range("B2").Value = "A"
$ws.Range("B3").Value = "B"

# Update the active selection to B3 (matches selection change in diff)
$ws.Range("B3").Select()
